$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.25625521910831139
$ws.Range("A2").Value = -0.063239506466892692
$ws.Range("A3").Value = -0.0039999999636091133
$ws.Range("A4").Value = -0.007999999934060753
$ws.Range("A5").Value = -0.0029999999597896121
$ws.Range("A6").Value = 0.0072863689492663752
$ws.Range("A7").Value = -0.009999999901460388
$ws.Range("A8").Value = -0.009999999900509593
$ws.Range("A9").Value = -0.0019999999496116416
$ws.Range("A10").Value = -0.0019999999493709453
$ws.Range("A11").Value = -0.0029999999431638003
$ws.Range("A12").Value = 0.052033956047663033
$ws.Range("A13").Value = -0.029760855793171004
$ws.Range("A14").Value = -0.0079999999011786116
$ws.Range("A15").Value = -0.00099999994393140668
$ws.Range("A16").Value = -0.0019999999370665655
$ws.Range("A17").Value = -0.0019999999360198473
$ws.Range("A18").Value = -0.0039999999233764072
$ws.Range("A19").Value = -0.0039999999728408397
$ws.Range("A20").Value = -0.0039999999707251988
$ws.Range("A21").Value = -0.0039999999704312117
$ws.Range("A22").Value = -0.0039999999701967326
$ws.Range("A23").Value = -0.0049999999556948893
$ws.Range("A24").Value = -0.019999999858495876
$ws.Range("A25").Value = -0.019999999856548101
$ws.Range("A26").Value = -0.0024999999466732703
$ws.Range("A27").Value = -0.0024999999442867349
$ws.Range("A28").Value = -0.0019999999339290753
$ws.Range("A29").Value = -0.006999999894181208
$ws.Range("A30").Value = -0.059999999560501394
$ws.Range("A31").Value = -0.0069999998860001966
$ws.Range("A32").Value = -0.0099999998668973689
$ws.Range("A33").Value = -0.0039999999039519452
